$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before row 1127, shifting the existing
# data (rows 1127:1187) down to 1129:1189 — mirrors the new weekly
# Acelga price-report rows being prepended to this block of the sheet.
$ws.Rows("1127:1128").Insert(-4121)

# Fill the two newly inserted rows with the new report data.

# Row 1127 - Calidad "Primera"
$ws.Range("A1127").Value = 6
$ws.Range("B1127").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1127").Value = "Metropolitana"
$ws.Range("D1127").Value = 45041
$ws.Range("E1127").Value = 13
$ws.Range("F1127").Value = 100112009
$ws.Range("G1127").Value = "Acelga"
$ws.Range("H1127").Value = "Sin especificar"
$ws.Range("I1127").Value = "Primera"
$ws.Range("J1127").Value = 470
$ws.Range("K1127").Value = 14000
$ws.Range("L1127").Value = 15000
$ws.Range("M1127").Value = 14468
$ws.Range("N1127").Value = "`$/docena de atados"
$ws.Range("O1127").Value = "Región Metropolitana"
$ws.Range("P1127").Value = 4823
$ws.Range("Q1127").Value = 3
$ws.Range("R1127").Value = "Hortaliza"

# Row 1128 - Calidad "Segunda"
$ws.Range("A1128").Value = 6
$ws.Range("B1128").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1128").Value = "Metropolitana"
$ws.Range("D1128").Value = 45041
$ws.Range("E1128").Value = 13
$ws.Range("F1128").Value = 100112009
$ws.Range("G1128").Value = "Acelga"
$ws.Range("H1128").Value = "Sin especificar"
$ws.Range("I1128").Value = "Segunda"
$ws.Range("J1128").Value = 170
$ws.Range("K1128").Value = 12000
$ws.Range("L1128").Value = 12000
$ws.Range("M1128").Value = 12000
$ws.Range("N1128").Value = "`$/docena de atados"
$ws.Range("O1128").Value = "Región Metropolitana"
$ws.Range("P1128").Value = 4000
$ws.Range("Q1128").Value = 3
$ws.Range("R1128").Value = "Hortaliza"
